# Adds two new "Title and Content" slides to the end of the deck:
#   Slide 10: "Java Beans"
#   Slide 11: "JSP and Servlet"
# mirroring the layout/structure already used by slides 2-9 (slideLayout2 /
# "Title and Content").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 - "Java Beans"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Java Beans"

$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "The beans of JavaBeans are classes that encapsulate one or more objects into a single standardized object."
$null = $body10.InsertAfter("`rIt is a reusable software component and can be manipulated visually in a builder tool (any IDE)..")
$null = $body10.InsertAfter("`rProperties:")
$null = $body10.InsertAfter("`rShould have a no-")
$null = $body10.InsertAfter("args")
$null = $body10.InsertAfter(" constructor")
$null = $body10.InsertAfter("`rHave getters and setters, private properties")
$null = $body10.InsertAfter("`rMust be serializable")

$body10.Paragraphs(4).IndentLevel = 2
$body10.Paragraphs(5).IndentLevel = 2
$body10.Paragraphs(6).IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 11 - "JSP and Servlet"
# ---------------------------------------------------------------------
$s11 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s11.Shapes.Item(1).TextFrame.TextRange.Text = "JSP and Servlet"

$body11 = $s11.Shapes.Item(2).TextFrame.TextRange
$body11.Text = "Client requests server for information and Server (Tomcat,  "
$null = $body11.InsertAfter("Jboss")
$null = $body11.InsertAfter(", ")
$null = $body11.InsertAfter("GlassFish")
$null = $body11.InsertAfter(") accesses the Web container in an application to process the request and provide response to the client.")
$null = $body11.InsertAfter("`rDeployment descriptor: web.xml ")
$null = $body11.InsertAfter(" Which request should call ")
$null = $body11.InsertAfter("which method")
$null = $body11.InsertAfter("`rThis contains servlets and servlet mappings")

$body11.Paragraphs(3).IndentLevel = 2

Write-Host "Added slides 10 and 11"
